$d = $word.ActiveDocument

# The second table on the page is the pairing-session log (Date | Time |
# Duration (h) | Driver | Observer). Its last existing row has its date
# ("02/03/2023") and time ("14:30") text split across multiple runs; use
# Find/Replace so the text collapses into a single run each (matching how
# Word normally stores freshly (re)typed text), leaving the values
# themselves unchanged.
$d.Content.Find.Execute("02/03/2023", $false, $false, $false, $false, $false, $true, 1, $false, "02/03/2023", 2) | Out-Null
$d.Content.Find.Execute("14:30", $false, $false, $false, $false, $false, $true, 1, $false, "14:30", 2) | Out-Null

# Locate the session-log table (the one whose header starts with "Date").
$logTable = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Cell(1, 1).Range.Text.TrimEnd([char]13, [char]7) -eq "Date") {
        $logTable = $tbl
    }
}

# Append two new pairing sessions recorded on 08/03/2023.
$row1 = $logTable.Rows.Add()
$row1.Cells.Item(1).Range.Text = "08/03/2023"
$row1.Cells.Item(2).Range.Text = "10:30"
$row1.Cells.Item(3).Range.Text = "1"
$row1.Cells.Item(4).Range.Text = "Driver"
$row1.Cells.Item(5).Range.Text = "Observer"

$row2 = $logTable.Rows.Add()
$row2.Cells.Item(1).Range.Text = "08/03/2023"
$row2.Cells.Item(2).Range.Text = "14:05"
$row2.Cells.Item(3).Range.Text = "3"
$row2.Cells.Item(4).Range.Text = "Observer"
$row2.Cells.Item(5).Range.Text = "Driver"
